# Update column G ("K" - strikeouts) with regenerated values.
# The author's commit regenerated save_data to use K (strikeouts) instead of
# the previous "Strike#" stat, recalculated std/mean, and wrote new s_vals.
# Net effect on this sheet: new numeric values in column G, rows 2-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 3
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 3
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 1
    26 = 0
    27 = 2
    28 = 2
    29 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
